$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: bold the "EX13.3 - Assumes use of STM32VL" header cell ---
$ws.Range("A5").Font.Bold = $true

# --- Row 6: rework the Exercise 13.2/13.3 prescaler calc ---
# B6 becomes a "Comma" styled cell holding the new prescaler_div value
$ws.Range("B6").Value = 4000000
$ws.Range("B6").Style = "Comma"

# C6 becomes a formula (=A6/B6) with a direct Comma number format applied
$ws.Range("C6").Formula = "=A6/B6"
$ws.Range("C6").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# F6 becomes a plain literal value (was a formula =408/3)
$ws.Range("F6").Value = 91

# --- Row 7: new "44100 Delta" label + delta formula ---
$ws.Range("F7").Value = "44100 Delta "
$ws.Range("F7").Font.Bold = $true
$ws.Range("G7").Formula = "=44100-G6"

# --- Remove the old scratch row 12 ---
$ws.Rows(12).Delete()

Write-Host "done"
